$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 9 - update title & link for pabii blog entry
$ws.Range("D9").Value = "해외대학이 대기업 취직에 유리한 이유 (1)"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/foreign-degree-job-market-merits-1/#utm_source=rss&utm_medium=rss&utm_campaign=foreign-degree-job-market-merits-1"

# Row 40 - update title for insightCampus entry
$ws.Range("D40").Value = "insightCampus"

# Row 46 - update title & link for BioinformaticsAndMe entry
$ws.Range("D46").Value = "[LG화학] 2021년 04월, 생물정보학(Bioinformatics 채용), 유전체 분석 및 DX 경력사원"
$ws.Range("E46").Value = "https://bioinformaticsandme.tistory.com/390"

# Row 50 - update title for incredible.egloos.com entry
$ws.Range("D50").Value = "인공신경망을 활용한 데이터 기반 결정구조 탐색"

# Row 51 - update title & link for bskyvision entry
$ws.Range("D51").Value = "[python] 리스트 안의 숫자 요소들을 문자로 변경하려면, map() 함수 사용"
$ws.Range("E51").Value = "https://bskyvision.com/1168"
